# Daily attendance processing - 2025-10-23 15:44:53
# Normalize/reorder the "Recorded By" (column G) contributor lists.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact-match replacement table for the "Recorded By" column values
# (old combined value -> new reordered value).
$map = @{
    "System, backup@backdoor.com, system" = "system, System, backup@backdoor.com";
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System";
    "System, admin@admin.com"             = "admin@admin.com, System";
    "admin@admin.com, dnasr281@gmail.com" = "dnasr281@gmail.com, admin@admin.com";
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    $val = $cell.Value2
    if ($val -ne $null -and $map.ContainsKey($val)) {
        $cell.Value2 = $map[$val]
    }
}
